$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Sheet1"

# Append the new water filling record in row 2
$ws.Range("A2").Value = "a"
$ws.Range("B2").Value = "b"
$ws.Range("C2").Value = "c"
$ws.Range("D2").Value = "d"
$ws.Range("E2").Value = "w"
$ws.Range("F2").Value = "e"
$ws.Range("G2").Value = "f"
$ws.Range("H2").Value = "g"
